$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Insert 4 new blank rows at row 2, pushing the existing data (old rows 2-21)
# down to rows 6-25.
$ws.Rows("2:5").Insert()

# Fill in the newly scraped records (match codes 4235, 4533, 4535, 4536).
# Columns A, C, D, E are stored as plain text; column B as a genuine number.
# Leading apostrophes force the numeric/percent-looking strings in A/C/D/E
# to be kept as literal text instead of being auto-detected as a number or
# percentage by Excel.

# Row 2: MATCH_CODE 4235 - no batting stats yet
$ws.Range("A2").Value = "'4235"
$ws.Range("F2").Value = "NO"

# Row 3: MATCH_CODE 4533 - no batting stats yet
$ws.Range("A3").Value = "'4533"
$ws.Range("F3").Value = "NO"

# Row 4: MATCH_CODE 4535
$ws.Range("A4").Value = "'4535"
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = "'0"
$ws.Range("D4").Value = "'0"
$ws.Range("E4").Value = "'1.27%"
$ws.Range("F4").Value = "NO"

# Row 5: MATCH_CODE 4536
$ws.Range("A5").Value = "'4536"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "'1"
$ws.Range("D5").Value = "'0"
$ws.Range("E5").Value = "'1.51%"
$ws.Range("F5").Value = "NO"

# Row insert in Excel copies formatting from the row above (the bold header
# row here), and the quoted literal-text entries above add a "quote prefix"
# style. Strip all of that back to the plain/default style used by the rest
# of the data rows, same as before the edit.
$ws.Range("A2:F5").ClearFormats()
